$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.027.22"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.541.09"
$ws.Range("E3").Value = "  -0.25%  "
$c = $ws.Range("D4")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$c = $ws.Range("D5")
$c.Value = "'603.54"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.34%  "
$c = $ws.Range("D6")
$c.Value = "'196.84"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.15%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -3.28%  "
$c = $ws.Range("D10")
$c.Value = "'0.652"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.12%  "
$c = $ws.Range("D11")
$c.Value = "'54.01"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").Value = "4.104.22"
$ws.Range("E14").Value = "  -0.16%  "
$c = $ws.Range("D15")
$c.Value = "'603.48"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "70.137.25"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "3.527.28"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  -0.56%  "
$c = $ws.Range("D22")
$c.Value = "'17.95"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.38%  "
$c = $ws.Range("D23")
$c.Value = "'5.30"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.61%  "
$c = $ws.Range("D24")
$c.Value = "'102.50"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("E25").Value = "  -3.01%  "
$c = $ws.Range("D26")
$c.Value = "'3.11"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.40%  "
$c = $ws.Range("D27")
$c.Value = "'10.97"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$c = $ws.Range("D28")
$c.Value = "'9.60"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("E29").Value = "  -1.84%  "
$c = $ws.Range("D30")
$c.Value = "'7.13"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "
$c = $ws.Range("D31")
$c.Value = "'4.31"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +16.97%  "
$c = $ws.Range("D32")
$c.Value = "'12.67"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  -1.86%  "
$c = $ws.Range("D34")
$c.Value = "'63.35"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "0.0₃0843"
$ws.Range("E35").Value = "  +7.60%  "
$ws.Range("D36").Value = "3.781.27"
$ws.Range("E36").Value = "  +6.92%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -3.27%  "
$c = $ws.Range("D39")
$c.Value = "'3.64"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  -1.90%  "
$c = $ws.Range("D41")
$c.Value = "'36.66"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.89%  "
$c = $ws.Range("D42")
$c.Value = "'495.94"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -7.99%  "
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("E46").Value = "  -2.18%  "
$c = $ws.Range("D47")
$c.Value = "'3.31"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("E50").Value = "  +2.14%  "
$c = $ws.Range("D51")
$c.Value = "'130.05"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.75%  "
